# Auto-generated Excel COM-interop script
# Applies numeric updates to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32: Automata for the People / Crab Oil
$ws.Range("H32").Value = 17424992
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 17424992
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 17424992
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -17425644

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 4328
$ws.Range("I40").Value = 1054.2858
$ws.Range("J40").Value = 7192.5
$ws.Range("K40").Value = 1054.2858
$ws.Range("L40").Value = 7192.5
$ws.Range("M40").Value = -879.2858000000001
$ws.Range("N40").Value = -7542.5

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 1058.12
$ws.Range("I98").Value = 1185.7
$ws.Range("J98").Value = 547.8
$ws.Range("K98").Value = 1185.7
$ws.Range("L98").Value = 547.8
$ws.Range("M98").Value = 312.3
$ws.Range("N98").Value = -3543.8

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 127357.94
$ws.Range("I116").Value = 178257.08
$ws.Range("J116").Value = 5200
$ws.Range("K116").Value = 178257.08
$ws.Range("L116").Value = 5200
$ws.Range("M116").Value = -174815.08
$ws.Range("N116").Value = -12084

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 1058.12
$ws.Range("I122").Value = 1185.7
$ws.Range("J122").Value = 547.8
$ws.Range("K122").Value = 3557.1
$ws.Range("L122").Value = 1643.4
$ws.Range("M122").Value = -1107.1
$ws.Range("N122").Value = -6543.4

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 794.44446
$ws.Range("I2").Value = 723
$ws.Range("J2").Value = 937.3333
$ws.Range("K2").Value = 723
$ws.Range("L2").Value = 937.3333
$ws.Range("M2").Value = -610
$ws.Range("N2").Value = -1163.3333

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 3097.48
$ws.Range("I32").Value = 2631.042
$ws.Range("J32").Value = 11959.8
$ws.Range("K32").Value = 2631.042
$ws.Range("L32").Value = 11959.8
$ws.Range("M32").Value = -2344.042
$ws.Range("N32").Value = -12533.8

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 6279.5713
$ws.Range("I63").Value = 7351.615
$ws.Range("K63").Value = 7351.615
$ws.Range("M63").Value = -6665.615

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 6279.5713
$ws.Range("I66").Value = 7351.615
$ws.Range("K66").Value = 36758.075
$ws.Range("M66").Value = -33326.075

# Row 82: Belle of the Brawl / Titanium Vambraces of Fending
$ws.Range("H82").Value = 27120
$ws.Range("J82").Value = 27120
$ws.Range("L82").Value = 27120
$ws.Range("N82").Value = -27842

# Row 85: Shouldering the Shut-ins (L) / Titanium Vambraces of Fending
$ws.Range("H85").Value = 27120
$ws.Range("J85").Value = 27120
$ws.Range("L85").Value = 27120
$ws.Range("N85").Value = -29616

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 1680.7727
$ws.Range("I110").Value = 1556.6842
$ws.Range("K110").Value = 1556.6842
$ws.Range("M110").Value = 488.3158000000001

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 794.44446
$ws.Range("I116").Value = 723
$ws.Range("J116").Value = 937.3333
$ws.Range("K116").Value = 723
$ws.Range("L116").Value = 937.3333
$ws.Range("M116").Value = 1571
$ws.Range("N116").Value = -5525.3333

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1976.6744
$ws.Range("I132").Value = 1676.8438
$ws.Range("J132").Value = 2848.9092
$ws.Range("K132").Value = 5030.5314
$ws.Range("L132").Value = 8546.7276
$ws.Range("M132").Value = -2500.5314
$ws.Range("N132").Value = -13606.7276

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 794.44446
$ws.Range("I3").Value = 723
$ws.Range("J3").Value = 937.3333
$ws.Range("K3").Value = 723
$ws.Range("L3").Value = 937.3333
$ws.Range("M3").Value = -609
$ws.Range("N3").Value = -1165.3333

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 7696207
$ws.Range("I99").Value = 2751959.8
$ws.Range("J99").Value = 25001072
$ws.Range("K99").Value = 2751959.8
$ws.Range("L99").Value = 25001072
$ws.Range("M99").Value = -2750461.8
$ws.Range("N99").Value = -25004068

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 23363.166
$ws.Range("I134").Value = 28910.611
$ws.Range("J134").Value = 6720.8335
$ws.Range("K134").Value = 86731.833
$ws.Range("L134").Value = 20162.5005
$ws.Range("M134").Value = -84196.833
$ws.Range("N134").Value = -25232.5005

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 2350.4
$ws.Range("I58").Value = 1695.4286
$ws.Range("J58").Value = 2787.0476
$ws.Range("K58").Value = 1695.4286
$ws.Range("L58").Value = 2787.0476
$ws.Range("M58").Value = -1492.4286
$ws.Range("N58").Value = -3193.0476

# Row 59: Bow Down to Magic / Crab Bow
$ws.Range("H59").Value = 64000
$ws.Range("J59").Value = 70000
$ws.Range("L59").Value = 70000
$ws.Range("N59").Value = -72290

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1232.919
$ws.Range("I134").Value = 930.96295
$ws.Range("K134").Value = 2792.88885
$ws.Range("M134").Value = -257.8888499999998

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 2350.4
$ws.Range("I136").Value = 1695.4286
$ws.Range("J136").Value = 2787.0476
$ws.Range("K136").Value = 5086.2858
$ws.Range("L136").Value = 8361.1428
$ws.Range("M136").Value = -2536.2858
$ws.Range("N136").Value = -13461.1428

$ws = $wb.Worksheets.Item("CUL")
# Row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 3990.182
$ws.Range("I137").Value = 2860
$ws.Range("J137").Value = 5346.4
$ws.Range("K137").Value = 8580
$ws.Range("L137").Value = 16039.2
$ws.Range("M137").Value = -3480
$ws.Range("N137").Value = -26239.2

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 8104.7666
$ws.Range("I70").Value = 3965
$ws.Range("K70").Value = 3965
$ws.Range("M70").Value = -3695

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 8104.7666
$ws.Range("I73").Value = 3965
$ws.Range("K73").Value = 3965
$ws.Range("M73").Value = -3029

# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 8801
$ws.Range("I107").Value = 11400.667
$ws.Range("K107").Value = 11400.667
$ws.Range("M107").Value = -9480.666999999999

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1421.6154
$ws.Range("I122").Value = 1228.7693
$ws.Range("J122").Value = 1614.4615
$ws.Range("K122").Value = 3686.3079
$ws.Range("L122").Value = 4843.3845
$ws.Range("M122").Value = -1236.3079
$ws.Range("N122").Value = -9743.3845

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore / Hard Leather
$ws.Range("H16").Value = 4097.933
$ws.Range("I16").Value = 3689.923
$ws.Range("J16").Value = 6750
$ws.Range("K16").Value = 3689.923
$ws.Range("L16").Value = 6750
$ws.Range("M16").Value = -3519.923
$ws.Range("N16").Value = -7090

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 643
$ws.Range("I46").Value = 678.8889
$ws.Range("J46").Value = 320
$ws.Range("K46").Value = 678.8889
$ws.Range("L46").Value = 320
$ws.Range("M46").Value = -490.8889
$ws.Range("N46").Value = -696

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 41668532
$ws.Range("I68").Value = 66668044
$ws.Range("J68").Value = 2677.7778
$ws.Range("K68").Value = 66668044
$ws.Range("L68").Value = 2677.7778
$ws.Range("M68").Value = -66667295
$ws.Range("N68").Value = -4175.7778

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 41668532
$ws.Range("I71").Value = 66668044
$ws.Range("J71").Value = 2677.7778
$ws.Range("K71").Value = 333340220
$ws.Range("L71").Value = 13388.889
$ws.Range("M71").Value = -333336476
$ws.Range("N71").Value = -20876.889

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 870
$ws.Range("I93").Value = 876.6667
$ws.Range("J93").Value = 850
$ws.Range("K93").Value = 876.6667
$ws.Range("L93").Value = 850
$ws.Range("M93").Value = 371.3333
$ws.Range("N93").Value = -3346

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 47000.184
$ws.Range("I122").Value = 56972.5
$ws.Range("K122").Value = 170917.5
$ws.Range("M122").Value = -168467.5

